$wb = $excel.ActiveWorkbook

# --- Macro_taxonomy sheet: update macro-proportion breakdown ---
$ws = $wb.Worksheets.Item("Macro_taxonomy")

# Split the Urban/Other row (was 100% MATO) into two rows: MATO 0.5 and ME+MEO/LWAL 0.5
$ws.Rows("16:16").Insert()
$ws.Range("D15").Value = 0.5
$ws.Range("A16").Value = "Other"
$ws.Range("B16").Value = "Urban"
$ws.Range("C16").Value = "ME+MEO/LWAL"
$ws.Range("D16").Value = 0.5

# Split the Rural/Other row (was 100% MATO) into three rows: MATO 0.5, EWV/LN 0.25, ME+MEO/LWAL 0.25
$ws.Range("D25").Value = 0.5
$ws.Range("A26").Value = "Other"
$ws.Range("B26").Value = "Rural"
$ws.Range("C26").Value = "EWV/LN"
$ws.Range("D26").Value = 0.25
$ws.Range("A27").Value = "Other"
$ws.Range("B27").Value = "Rural"
$ws.Range("C27").Value = "ME+MEO/LWAL"
$ws.Range("D27").Value = 0.25

$ws.Range("D17").Select() | Out-Null
$ws.Activate() | Out-Null

# --- Costs sheet loses the "active tab" flag as Macro_taxonomy becomes the active sheet ---
$wb.Worksheets.Item("Costs").Range("E1:E33").Select() | Out-Null
$ws.Activate() | Out-Null
